$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 177.72728
$ws.Range("J41").Value = 97
$ws.Range("L41").Value = 97
$ws.Range("N41").Value = -977

# Row 62
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

# Row 65
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

# Row 70
$ws.Range("H70").Value = 6187.1875
$ws.Range("I70").Value = 5999.6665
$ws.Range("J70").Value = 6749.75
$ws.Range("K70").Value = 17998.9995
$ws.Range("L70").Value = 20249.25
$ws.Range("M70").Value = -17728.9995
$ws.Range("N70").Value = -20789.25

# Row 73
$ws.Range("H73").Value = 6187.1875
$ws.Range("I73").Value = 5999.6665
$ws.Range("J73").Value = 6749.75
$ws.Range("K73").Value = 17998.9995
$ws.Range("L73").Value = 20249.25
$ws.Range("M73").Value = -17062.9995
$ws.Range("N73").Value = -22121.25

# Row 96
$ws.Range("H96").Value = 1884.3334
$ws.Range("I96").Value = 2371.75
$ws.Range("J96").Value = 909.5
$ws.Range("K96").Value = 7115.25
$ws.Range("L96").Value = 2728.5
$ws.Range("M96").Value = -5742.25
$ws.Range("N96").Value = -5474.5

# Row 98
$ws.Range("H98").Value = 1238.5555
$ws.Range("I98").Value = 1193.5
$ws.Range("K98").Value = 1193.5
$ws.Range("M98").Value = 304.5

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents() | Out-Null
$ws.Range("N113").ClearContents() | Out-Null

# Row 122
$ws.Range("H122").Value = 1238.5555
$ws.Range("I122").Value = 1193.5
$ws.Range("K122").Value = 3580.5
$ws.Range("M122").Value = -1130.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2474.25
$ws.Range("I2").Value = 1500
$ws.Range("K2").Value = 1500
$ws.Range("M2").Value = -1387

# Row 45
$ws.Range("H45").Value = 2125
$ws.Range("I45").Value = 1542.5
$ws.Range("J45").Value = 4066.6667
$ws.Range("K45").Value = 1542.5
$ws.Range("L45").Value = 4066.6667
$ws.Range("M45").Value = -1165.5
$ws.Range("N45").Value = -4820.6667

# Row 61
$ws.Range("H61").Value = 7013.4165
$ws.Range("I61").Value = 7196.4546
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 7196.4546
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -6984.4546
$ws.Range("N61").Value = -5424

# Row 74
$ws.Range("H74").Value = 1511.25
$ws.Range("I74").Value = 926.6667
$ws.Range("K74").Value = 926.6667
$ws.Range("M74").Value = -52.66669999999999

# Row 77
$ws.Range("H77").Value = 1511.25
$ws.Range("I77").Value = 926.6667
$ws.Range("K77").Value = 4633.3335
$ws.Range("M77").Value = -265.3334999999997

# Row 102
$ws.Range("H102").Value = 7940628
$ws.Range("I102").Value = 8551292
$ws.Range("K102").Value = 8551292
$ws.Range("M102").Value = -8549670

# Row 116
$ws.Range("H116").Value = 2474.25
$ws.Range("I116").Value = 1500
$ws.Range("K116").Value = 1500
$ws.Range("M116").Value = 794

# Row 136
$ws.Range("H136").Value = 7013.4165
$ws.Range("I136").Value = 7196.4546
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 21589.3638
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -19039.3638
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2474.25
$ws.Range("I3").Value = 1500
$ws.Range("K3").Value = 1500
$ws.Range("M3").Value = -1386

# Row 86
$ws.Range("H86").Value = 1087.125
$ws.Range("I86").Value = 924.25
$ws.Range("K86").Value = 924.25
$ws.Range("M86").Value = 198.75

# Row 89
$ws.Range("H89").Value = 1087.125
$ws.Range("I89").Value = 924.25
$ws.Range("K89").Value = 4621.25
$ws.Range("M89").Value = 994.75

# Row 134
$ws.Range("H134").Value = 3400
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents() | Out-Null

# Row 138
$ws.Range("H138").Value = 59999
$ws.Range("I138").Value = 59999
$ws.Range("K138").Value = 59999
$ws.Range("M138").Value = -54859

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 16170.318
$ws.Range("I99").Value = 11222.5
$ws.Range("K99").Value = 11222.5
$ws.Range("M99").Value = -9724.5

# Row 122
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents() | Out-Null

# Row 126
$ws.Range("H126").Value = 16170.318
$ws.Range("I126").Value = 11222.5
$ws.Range("K126").Value = 33667.5
$ws.Range("M126").Value = -31197.5

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents() | Out-Null

# Row 29
$ws.Range("H29").Value = 17348.4
$ws.Range("I29").Value = 19185.5
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 19185.5
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -18895.5
$ws.Range("N29").Value = -10580

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2005
$ws.Range("J7").Value = 2005
$ws.Range("L7").Value = 2005
$ws.Range("N7").Value = -2229

# Row 33
$ws.Range("H33").Value = 5049666.5
$ws.Range("I33").Value = 7537000
$ws.Range("J33").Value = 75000
$ws.Range("K33").Value = 7537000
$ws.Range("L33").Value = 75000
$ws.Range("M33").Value = -7536710
$ws.Range("N33").Value = -75580

# Row 40
$ws.Range("H40").Value = 3135.7917
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents() | Out-Null

# Row 47
$ws.Range("H47").Value = 26250
$ws.Range("J47").Value = 26250
$ws.Range("L47").Value = 26250
$ws.Range("N47").Value = -27230

# Row 52
$ws.Range("H52").Value = 26250
$ws.Range("J52").Value = 26250
$ws.Range("L52").Value = 26250
$ws.Range("N52").Value = -26716

# Row 100
$ws.Range("H100").Value = 2688.125
$ws.Range("I100").Value = 2600.7144
$ws.Range("K100").Value = 2600.7144
$ws.Range("M100").Value = -2059.7144

# Row 126
$ws.Range("H126").Value = 2005
$ws.Range("J126").Value = 2005
$ws.Range("L126").Value = 6015
$ws.Range("N126").Value = -10955

# Row 132
$ws.Range("H132").Value = 1599.7142
$ws.Range("I132").Value = 1399.6364
$ws.Range("J132").Value = 2333.3333
$ws.Range("K132").Value = 4198.9092
$ws.Range("L132").Value = 6999.999899999999
$ws.Range("M132").Value = -1668.9092
$ws.Range("N132").Value = -12059.9999

# Row 136
$ws.Range("H136").Value = 7334
$ws.Range("I136").Value = 7334
$ws.Range("K136").Value = 22002
$ws.Range("M136").Value = -19452

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2112.5
$ws.Range("I81").Value = 2112.5
$ws.Range("K81").Value = 4225
$ws.Range("M81").Value = -3164

# Row 84
$ws.Range("H84").Value = 2112.5
$ws.Range("I84").Value = 2112.5
$ws.Range("K84").Value = 21125
$ws.Range("M84").Value = -15821
